# Updates cryptos price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.790.82'
$ws.Range("E2").Value = '  +1.59%  '
$ws.Range("D3").Value = '3.116.93'
$ws.Range("E3").Value = '  +1.27%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'532.65"
$ws.Range("E5").Value = '  +2.43%  '
$ws.Range("D6").Value = "'138.25"
$ws.Range("E6").Value = '  +1.63%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").Value = "'0.496"
$ws.Range("E8").Value = '  +10.05%  '
$ws.Range("D9").Value = "'7.36"
$ws.Range("E9").Value = '  +0.53%  '
$ws.Range("E10").Value = '  +1.47%  '
$ws.Range("E11").Value = '  +3.95%  '
$ws.Range("E12").Value = '  +3.55%  '
$ws.Range("D13").Value = '3.657.47'
$ws.Range("E13").Value = '  +1.26%  '
$ws.Range("D14").Value = "'25.77"
$ws.Range("E14").Value = '  +1.66%  '
$ws.Range("E15").Value = '  +3.16%  '
$ws.Range("D16").Value = '57.895.39'
$ws.Range("E16").Value = '  +1.60%  '
$ws.Range("D17").Value = '3.122.77'
$ws.Range("E17").Value = '  +1.18%  '
$ws.Range("D18").Value = "'6.13"
$ws.Range("E18").Value = '  +4.55%  '
$ws.Range("D19").Value = "'12.83"
$ws.Range("E19").Value = '  +3.12%  '
$ws.Range("D20").Value = "'8.11"
$ws.Range("E20").Value = '  +3.43%  '
$ws.Range("D21").Value = "'373.10"
$ws.Range("E21").Value = '  +7.68%  '
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").Value = "'5.72"
$ws.Range("E23").Value = '  -1.34%  '
$ws.Range("D24").Value = "'69.30"
$ws.Range("E24").Value = '  +1.54%  '
$ws.Range("E25").Value = '  +2.20%  '
$ws.Range("E26").Value = '  -0.21%  '
$ws.Range("E27").Value = '  +0.34%  '
$ws.Range("D28").Value = '0.0₃0878'
$ws.Range("E28").Value = '  +0.17%  '
$ws.Range("E29").Value = '  +4.22%  '
$ws.Range("D30").Value = "'6.15"
$ws.Range("E30").Value = '  +4.49%  '
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("D32").Value = "'21.52"
$ws.Range("E32").Value = '  +3.72%  '
$ws.Range("D33").Value = "'5.16"
$ws.Range("E33").Value = '  +5.18%  '
$ws.Range("E34").Value = '  +3.11%  '
$ws.Range("D35").Value = "'160.27"
$ws.Range("E35").Value = '  +0.79%  '
$ws.Range("E36").Value = '  +2.93%  '
$ws.Range("D37").Value = "'1.30"
$ws.Range("E37").Value = '  +5.06%  '
$ws.Range("D38").Value = "'25.61"
$ws.Range("E38").Value = '  -0.46%  '
$ws.Range("E39").Value = '  +3.72%  '
$ws.Range("D40").Value = "'0.0670"
$ws.Range("E40").Value = '  +2.88%  '
$ws.Range("D41").Value = '2.555.04'
$ws.Range("E41").Value = '  +7.48%  '
$ws.Range("D42").Value = "'4.14"
$ws.Range("E42").Value = '  +3.24%  '
$ws.Range("D43").Value = "'38.29"
$ws.Range("E43").Value = '  +4.76%  '
$ws.Range("D44").Value = "'0.696"
$ws.Range("E44").Value = '  +1.06%  '
$ws.Range("E45").Value = '  +3.41%  '
$ws.Range("E46").Value = '  -0.01%  '
$ws.Range("D47").Value = "'0.980"
$ws.Range("E47").Value = '  +2.06%  '
$ws.Range("D50").Value = "'0.0946"
$ws.Range("E50").Value = '  +6.39%  '
$ws.Range("D51").Value = "'0.746"
$ws.Range("E51").Value = '  -0.81%  '
